$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Q03" quiz column (D) ---------------------------------
$ws.Range("D1").Value = "Q03"

# Row 2-15: per-student Q03 scores (mirrors the Q01/Q02 "(score/total)*10"
# pattern already used in column C, with a couple of literal zeros).
$ws.Range("D2").Formula  = "=(12/20)*10"
$ws.Range("D3").Formula  = "=(7/20)*10"
$ws.Range("D4").Formula  = "=(12/20)*10"
$ws.Range("D5").Formula  = "=(10/20)*10"
$ws.Range("D6").Formula  = "=(11/20)*10"
$ws.Range("D7").Formula  = "=(16/20)*10"
$ws.Range("D8").Value    = 0
$ws.Range("D9").Value    = 0
$ws.Range("D10").Value   = 0
$ws.Range("D11").Formula = "=(13/20)*10"
$ws.Range("D12").Formula = "=(6/20)*10"
$ws.Range("D13").Formula = "=(13/20)*10"
$ws.Range("D14").Value   = 0
$ws.Range("D15").Value   = 0

# --- Update the active selection (author moved to C36 afterwards) -------
$ws.Range("C36").Select()
